$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B13 value (350 -> 253)
$ws.Range("B13").Value = "253"

# Add new row 14 content: invalid productcode negative test case
$ws.Range("A14").Value = "invalid.productcode"
$ws.Range("B14").Value = "test"

# Move selection to B15
$ws.Range("B15").Select()
